$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 33335344
$ws.Range("I32").Value = 1334.75
$ws.Range("K32").Value = 1334.75
$ws.Range("M32").Value = -1008.75

$ws.Range("H33").Value = 8640797
$ws.Range("I33").Value = 12626833
$ws.Range("J33").Value = 4383.8335
$ws.Range("K33").Value = 12626833
$ws.Range("L33").Value = 4383.8335
$ws.Range("M33").Value = -12626604
$ws.Range("N33").Value = -4841.8335

$ws.Range("H40").Value = 1148.9
$ws.Range("J40").Value = 1148.9
$ws.Range("L40").Value = 1148.9
$ws.Range("N40").Value = -1498.9

$ws.Range("H106").Value = 1236.5428
$ws.Range("I106").Value = 1314.963
$ws.Range("K106").Value = 1314.963
$ws.Range("M106").Value = -683.963

$ws.Range("H137").Value = 5556732.5
$ws.Range("I137").Value = 1246.4117
$ws.Range("K137").Value = 3739.2351
$ws.Range("M137").Value = -1189.2351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 601.3333
$ws.Range("I4").Value = 556.8182
$ws.Range("K4").Value = 556.8182
$ws.Range("M4").Value = -440.8182

$ws.Range("H45").Value = 79899.46000000001
$ws.Range("I45").Value = 103029.1
$ws.Range("K45").Value = 103029.1
$ws.Range("M45").Value = -102652.1

$ws.Range("H140").Value = 67500
$ws.Range("J140").Value = 67500
$ws.Range("L140").Value = 67500
$ws.Range("N140").Value = -77860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H86").Value = 1845.4736
$ws.Range("I86").Value = 1035.9375
$ws.Range("J86").Value = 6163
$ws.Range("K86").Value = 1035.9375
$ws.Range("L86").Value = 6163
$ws.Range("M86").Value = 87.0625
$ws.Range("N86").Value = -8409

$ws.Range("H89").Value = 1845.4736
$ws.Range("I89").Value = 1035.9375
$ws.Range("J89").Value = 6163
$ws.Range("K89").Value = 5179.6875
$ws.Range("L89").Value = 30815
$ws.Range("M89").Value = 436.3125
$ws.Range("N89").Value = -42047

$ws.Range("H107").Value = 7531.3145
$ws.Range("I107").Value = 9089.925999999999
$ws.Range("K107").Value = 9089.925999999999
$ws.Range("M107").Value = -7169.925999999999

$ws.Range("H123").Value = 187777
$ws.Range("J123").Value = 187777
$ws.Range("L123").Value = 187777
$ws.Range("N123").Value = -197577

$ws.Range("H132").Value = 66666.664
$ws.Range("I132").Value = 50000
$ws.Range("J132").Value = 150000
$ws.Range("K132").Value = 50000
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -44940
$ws.Range("N132").Value = -160120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 830.449
$ws.Range("I22").Value = 813.68085
$ws.Range("K22").Value = 813.68085
$ws.Range("M22").Value = -463.68085

$ws.Range("H58").Value = 1995.64
$ws.Range("I58").Value = 2406.3076
$ws.Range("J58").Value = 1550.75
$ws.Range("K58").Value = 2406.3076
$ws.Range("L58").Value = 1550.75
$ws.Range("M58").Value = -2203.3076
$ws.Range("N58").Value = -1956.75

$ws.Range("H106").Value = 52804
$ws.Range("J106").Value = 52804
$ws.Range("L106").Value = 52804
$ws.Range("N106").Value = -55328

$ws.Range("H107").Value = 1819.5161
$ws.Range("I107").Value = 1576.2858
$ws.Range("K107").Value = 1576.2858
$ws.Range("M107").Value = 343.7141999999999

$ws.Range("H134").Value = 1408.7931
$ws.Range("I134").Value = 1335.0741
$ws.Range("J134").Value = 2404
$ws.Range("K134").Value = 4005.2223
$ws.Range("L134").Value = 7212
$ws.Range("M134").Value = -1470.2223
$ws.Range("N134").Value = -12282

$ws.Range("H136").Value = 1995.64
$ws.Range("I136").Value = 2406.3076
$ws.Range("J136").Value = 1550.75
$ws.Range("K136").Value = 7218.9228
$ws.Range("L136").Value = 4652.25
$ws.Range("M136").Value = -4668.9228
$ws.Range("N136").Value = -9752.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3509.3333
$ws.Range("I34").Value = 348
$ws.Range("J34").Value = 3796.7273
$ws.Range("K34").Value = 1044
$ws.Range("L34").Value = 11390.1819
$ws.Range("M34").Value = -960
$ws.Range("N34").Value = -11558.1819

$ws.Range("H68").Value = 2148.2
$ws.Range("I68").Value = 1664.1666
$ws.Range("K68").Value = 4992.4998
$ws.Range("M68").Value = -4181.4998

$ws.Range("H71").Value = 2148.2
$ws.Range("I71").Value = 1664.1666
$ws.Range("K71").Value = 14977.4994
$ws.Range("M71").Value = -10921.4994

$ws.Range("H81").Value = 23815726
$ws.Range("J81").Value = 8092.5
$ws.Range("L81").Value = 24277.5
$ws.Range("N81").Value = -26523.5

$ws.Range("H84").Value = 23815726
$ws.Range("J84").Value = 8092.5
$ws.Range("L84").Value = 72832.5
$ws.Range("N84").Value = -84064.5

$ws.Range("H109").Value = 4844
$ws.Range("I109").Value = 2678.8572
$ws.Range("K109").Value = 8036.571599999999
$ws.Range("M109").Value = -6996.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 23926790
$ws.Range("I80").Value = 125365.78
$ws.Range("J80").Value = 41777856
$ws.Range("K80").Value = 125365.78
$ws.Range("L80").Value = 41777856
$ws.Range("M80").Value = -124367.78
$ws.Range("N80").Value = -41779852

$ws.Range("H83").Value = 23926790
$ws.Range("I83").Value = 125365.78
$ws.Range("J83").Value = 41777856
$ws.Range("K83").Value = 626828.9
$ws.Range("L83").Value = 208889280
$ws.Range("M83").Value = -621836.9
$ws.Range("N83").Value = -208899264

$ws.Range("H122").Value = 2793.8215
$ws.Range("I122").Value = 2843.4211
$ws.Range("J122").Value = 2689.111
$ws.Range("K122").Value = 8530.263300000001
$ws.Range("L122").Value = 8067.333
$ws.Range("M122").Value = -6080.263300000001
$ws.Range("N122").Value = -12967.333

$ws.Range("H136").Value = 35000
$ws.Range("J136").Value = 35000
$ws.Range("L136").Value = 105000
$ws.Range("N136").Value = -110100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 99999
$ws.Range("J20").Value = 99999
$ws.Range("L20").Value = 99999
$ws.Range("N20").Value = -100451

$ws.Range("H22").Value = 6111.5625
$ws.Range("I22").Value = 2833.3333
$ws.Range("J22").Value = 6868.077
$ws.Range("K22").Value = 2833.3333
$ws.Range("L22").Value = 6868.077
$ws.Range("M22").Value = -2538.3333
$ws.Range("N22").Value = -7458.077

$ws.Range("H27").Value = 6111.5625
$ws.Range("I27").Value = 2833.3333
$ws.Range("J27").Value = 6868.077
$ws.Range("K27").Value = 2833.3333
$ws.Range("L27").Value = 6868.077
$ws.Range("M27").Value = -2726.3333
$ws.Range("N27").Value = -7082.077

$ws.Range("H43").Value = 2477727.2
$ws.Range("J43").Value = 3390000
$ws.Range("L43").Value = 3390000
$ws.Range("N43").Value = -3390386

$ws.Range("H46").Value = 12525.333
$ws.Range("I46").Value = 18925.857
$ws.Range("J46").Value = 6924.875
$ws.Range("K46").Value = 18925.857
$ws.Range("L46").Value = 6924.875
$ws.Range("M46").Value = -18737.857
$ws.Range("N46").Value = -7300.875

$ws.Range("H68").Value = 3591.5715
$ws.Range("I68").Value = 3387
$ws.Range("J68").Value = 3864.3333
$ws.Range("K68").Value = 3387
$ws.Range("L68").Value = 3864.3333
$ws.Range("M68").Value = -2638
$ws.Range("N68").Value = -5362.3333

$ws.Range("H71").Value = 3591.5715
$ws.Range("I71").Value = 3387
$ws.Range("J71").Value = 3864.3333
$ws.Range("K71").Value = 16935
$ws.Range("L71").Value = 19321.6665
$ws.Range("M71").Value = -13191
$ws.Range("N71").Value = -26809.6665

$ws.Range("H122").Value = 2589.6667
$ws.Range("I122").Value = 2440.16
$ws.Range("K122").Value = 7320.48
$ws.Range("M122").Value = -4870.48

$ws.Range("H125").Value = 190000
$ws.Range("J125").Value = 190000
$ws.Range("L125").Value = 190000
$ws.Range("N125").Value = -199840

$ws.Range("H132").Value = 2398.9812
$ws.Range("I132").Value = 1982.081
$ws.Range("J132").Value = 3363.0625
$ws.Range("K132").Value = 5946.242999999999
$ws.Range("L132").Value = 10089.1875
$ws.Range("M132").Value = -3416.242999999999
$ws.Range("N132").Value = -15149.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2606.1738
$ws.Range("I126").Value = 2228.5
$ws.Range("J126").Value = 3018.182
$ws.Range("K126").Value = 6685.5
$ws.Range("L126").Value = 9054.545999999998
$ws.Range("M126").Value = -4215.5
$ws.Range("N126").Value = -13994.546
